$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "30.063.01"
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +5.47%  "

$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "1.921.48"
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +2.48%  "

$ws.Cells.Item(4, 5).Value = "  -0.80%  "

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "332.34"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +5.22%  "

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "1.0000"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -0.80%  "

$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "0.5251"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +3.23%  "

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "0.4055"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +3.90%  "

$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.08554"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +2.30%  "

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "43.01"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +3.45%  "

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "1.130"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +2.33%  "

$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "22.32"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +9.31%  "

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "6.416"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +3.00%  "

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "1.919.44"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +2.48%  "

$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "7.406"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +1.91%  "

$ws.Cells.Item(16, 5).Value = "  -0.83%  "

$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "96.57"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +5.81%  "

$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "0.00001117"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +1.13%  "

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "0.06704"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -0.36%  "

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "18.29"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +3.31%  "

$ws.Cells.Item(21, 5).Value = "  -0.75%  "

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "6.072"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +2.65%  "

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "30.076.63"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +5.44%  "

$ws.Cells.Item(24, 5).Value = "  +1.39%  "

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "2.227"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -0.07%  "

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "2.140.13"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +2.49%  "

$ws.Cells.Item(27, 5).Value = "  +2.65%  "

$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "160.83"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -0.82%  "

$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "2.459"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +2.63%  "

$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "129.63"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +3.23%  "

$ws.Cells.Item(31, 5).Value = "  +3.85%  "

$ws.Cells.Item(32, 5).Value = "  +1.58%  "

$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "6.116"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +6.18%  "

$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "3.650"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +1.04%  "

$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "0.02530"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +2.97%  "

$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "0.06618"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +1.44%  "

$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "0.2232"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +3.28%  "

$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "1.240"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +4.18%  "

$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "9.055"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +2.76%  "

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "5.224"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +3.51%  "

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "0.6574"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +2.95%  "

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "11.71"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +5.74%  "

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "1.246"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +0.57%  "

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "0.6215"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +3.61%  "

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "13.29"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +1.75%  "

$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "3.794"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +3.01%  "

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "2.089"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +4.13%  "

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "1.244"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +2.42%  "

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "125.74"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +3.20%  "

$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "80.32"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +5.22%  "

$ws.Cells.Item(51, 5).Value = "  +0.77%  "
